$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.288.98"
$ws.Range("E2").Value = "  -4.66%  "
$ws.Range("D3").Value = "2.238.92"
$ws.Range("E3").Value = "  -5.87%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.586"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.46%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0828"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "2.580.61"
$ws.Range("E15").Value = "  -5.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.867"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.20%  "
$ws.Range("D18").Value = "2.228.06"
$ws.Range("E18").Value = "  -5.59%  "
$ws.Range("D19").Value = "43.234.74"
$ws.Range("E19").Value = "  -4.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.14%  "
$ws.Range("D21").Value = "0.0₃0968"
$ws.Range("E21").Value = "  -9.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.64%  "
$ws.Range("E24").Value = "  -12.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "237.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -15.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.61%  "
$ws.Range("E34").Value = "  -8.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.22%  "
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  -8.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.87%  "
$ws.Range("E41").Value = "  -11.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0325"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "1.780.29"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.87%  "
$ws.Range("E48").Value = "  -10.40%  "
$ws.Range("E49").Value = "  -11.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.95%  "
